$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B15 was stored as text "4"; convert to a real number
$ws.Range("B15").Value = 4

# Add new row 16 with annotation data
$ws.Range("A16").Value = "Ruilin"
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = "3"
$ws.Range("C16").Value = "无"
$ws.Range("D16").Value = "DFT"
$ws.Range("E16").Value = "MET"
$ws.Range("F16").Value = "4cbdf296-0ef7-4a60-9d08-bf70fb941ab3"
$ws.Range("G16").Value = "SJTB5GZCb_annotated.xlsx"
$ws.Range("H16").Value = "The paper does not sufficiently discuss and compare the relevant neuroscience literature and related work."
